$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 23 (existing rows 23+ shift down to 26+)
$ws.Rows("23:25").Insert()

# Reset formatting on the freshly inserted rows to the sheet's default "wrap
# text" look (style used by every other data row), then re-apply the
# non-wrapping look used in column C for this block (copied from C22, the
# row immediately above, which already carries that exact style). Column H
# is left completely empty/unformatted (unused in this block, same as the
# surrounding rows).
$newRows = $ws.Range("A23:I25")
$newRows.ClearFormats()
$newRows.WrapText = $true
$ws.Range("H23:H25").Clear()

$ws.Range("C22").Copy()
$ws.Range("C23:C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 23: options_num_slices ---
$ws.Range("A23").Value = "options_num_slices"
$ws.Range("B23").Value = "Number of slices"
$ws.Range("C23").Value = "Number of slices"
$ws.Range("D23").Value = "number"
$ws.Range("E23").Value = "int"
$ws.Range("F23").Value = "Check with Vince"
$ws.Range("G23").Value = "Check with Vince"
$ws.Range("I23").Value = "slicetiming"
$ws.Rows(23).RowHeight = 17

# --- Row 24: options_repetition_time ---
$ws.Range("A24").Value = "options_repetition_time"
$ws.Range("B24").Value = "Repetition time(secs)"
$ws.Range("C24").Value = "Repetition time in secs"
$ws.Range("D24").Value = "number"
$ws.Range("E24").Value = "float"
$ws.Range("F24").Value = "Check with Vince"
$ws.Range("G24").Value = "Check with Vince"
$ws.Range("I24").Value = "slicetiming"
$ws.Rows(24).RowHeight = 17

# --- Row 25: options_acquisition_order ---
$ws.Range("A25").Value = "options_acquisition_order"
$ws.Range("B25").Value = "Acquisition order"
$ws.Range("C25").Value = "Slicetime acquisition order"
$ws.Range("D25").Value = "array of number items. Size of array MUST BE EQUAL TO number of slices. So you can do this check in the UI incase the customer enters the num_slices and acquisition_order"
$ws.Range("E25").Value = "array with int items. Size of array MUST BE EQUAL TO number of slices."
$ws.Range("F25").Value = "Check with Vince"
$ws.Range("G25").Value = "Check with Vince"
$ws.Range("I25").Value = "slicetiming"
$ws.Rows(25).RowHeight = 306

# Update sheet view to match the saved workbook state (selection moves to
# A26, the first cell of the block that used to be row 23).
$ws.Range("A26").Select()
